$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure these price cells stay text (they use "." as thousands separators,
# not decimal points), matching the source data which is plain text.
$textCells = @("D4", "D5", "D6", "D7", "D9", "D10", "D11", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D22", "D23", "D24", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range("D2").Value = "29.457.61"
$ws.Range("E2").Value = "  +0.53%  "
$ws.Range("D3").Value = "1.909.37"
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  +0.70%  "
$ws.Range("D5").Value = "325.37"
$ws.Range("E5").Value = "  +1.04%  "
$ws.Range("D6").Value = "1.007"
$ws.Range("E6").Value = "  +0.59%  "
$ws.Range("D7").Value = "0.4824"
$ws.Range("E7").Value = "  +2.40%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "0.08153"
$ws.Range("E9").Value = "  +1.42%  "
$ws.Range("D10").Value = "1.012"
$ws.Range("D11").Value = "23.43"
$ws.Range("E11").Value = "  +3.22%  "
$ws.Range("D12").Value = "1.929.07"
$ws.Range("E12").Value = "  -1.97%  "
$ws.Range("D13").Value = "6.013"
$ws.Range("E13").Value = "  +2.21%  "
$ws.Range("D14").Value = "7.163"
$ws.Range("E14").Value = "  +0.89%  "
$ws.Range("D15").Value = "90.33"
$ws.Range("E15").Value = "  +0.85%  "
$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D16").Value = "1.008"
$ws.Range("E16").Value = "  +0.68%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "0.06772"
$ws.Range("E17").Value = "  +2.20%  "
$ws.Range("D18").Value = "0.00001036"
$ws.Range("E18").Value = "  +0.61%  "
$ws.Range("D19").Value = "17.69"
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("D20").Value = "1.007"
$ws.Range("D21").Value = "29.468.72"
$ws.Range("E21").Value = "  +0.49%  "
$ws.Range("D22").Value = "5.629"
$ws.Range("E22").Value = "  +1.82%  "
$ws.Range("D23").Value = "11.72"
$ws.Range("E23").Value = "  +2.57%  "
$ws.Range("D24").Value = "2.177"
$ws.Range("E24").Value = "  -1.18%  "
$ws.Range("D25").Value = "2.141.48"
$ws.Range("E25").Value = "  -1.74%  "
$ws.Range("D26").Value = "156.39"
$ws.Range("E26").Value = "  +0.96%  "
$ws.Range("D27").Value = "6.479"
$ws.Range("E27").Value = "  +7.76%  "
$ws.Range("D28").Value = "20.05"
$ws.Range("E28").Value = "  +1.39%  "
$ws.Range("D29").Value = "2.109"
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").Value = "120.36"
$ws.Range("E30").Value = "  +2.17%  "
$ws.Range("D31").Value = "1.026"
$ws.Range("E31").Value = "  -4.19%  "
$ws.Range("D32").Value = "0.09524"
$ws.Range("E32").Value = "  +0.30%  "
$ws.Range("D33").Value = "5.514"
$ws.Range("E33").Value = "  +2.44%  "
$ws.Range("D34").Value = "3.564"
$ws.Range("E34").Value = "  +0.54%  "
$ws.Range("D35").Value = "1.390"
$ws.Range("E35").Value = "  -2.20%  "
$ws.Range("D36").Value = "0.02270"
$ws.Range("E36").Value = "  +0.91%  "
$ws.Range("D37").Value = "0.06104"
$ws.Range("E37").Value = "  +0.40%  "
$ws.Range("D39").Value = "10.82"
$ws.Range("E39").Value = "  +6.95%  "
$ws.Range("D40").Value = "0.5946"
$ws.Range("E40").Value = "  +1.61%  "
$ws.Range("D41").Value = "7.982"
$ws.Range("E41").Value = "  -3.02%  "
$ws.Range("D42").Value = "0.1855"
$ws.Range("E42").Value = "  +1.01%  "
$ws.Range("D43").Value = "1.277"
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("D44").Value = "2.373"
$ws.Range("E44").Value = "  -5.32%  "
$ws.Range("D45").Value = "12.56"
$ws.Range("E45").Value = "  +4.04%  "
$ws.Range("D46").Value = "0.07610"
$ws.Range("E46").Value = "  -3.83%  "
$ws.Range("D47").Value = "0.5565"
$ws.Range("E47").Value = "  +1.07%  "
$ws.Range("D48").Value = "1.946"
$ws.Range("E48").Value = "  +1.25%  "
$ws.Range("D49").Value = "116.58"
$ws.Range("E49").Value = "  +3.06%  "
$ws.Range("D50").Value = "72.52"
$ws.Range("E50").Value = "  +1.62%  "
$ws.Range("D51").Value = "2.407"
$ws.Range("E51").Value = "  +2.46%  "
